{"js": "// Replace each two-digit-by-two-digit multiplication problem/answer\n// with its updated version. Every \"old\" string below is unique in the\n// document, so a matchCase + matchWholeWord search locates exactly the\n// one run that needs to change.\nconst replacements = [\n  [\"25\u00d750=1250\", \"21\u00d750=1050\"],\n  [\"94\u00d739=3666\", \"44\u00d718=792\"],\n  [\"39\u00d753=2067\", \"64\u00d715=960\"],\n  [\"72\u00d781=5832\", \"61\u00d794=5734\"],\n  [\"12\u00d752=624\", \"44\u00d733=1452\"],\n  [\"45\u00d741=1845\", \"12\u00d750=600\"],\n  [\"64\u00d748=3072\", \"61\u00d743=2623\"],\n  [\"22\u00d769=1518\", \"45\u00d738=1710\"],\n  [\"76\u00d743=3268\", \"75\u00d788=6600\"],\n  [\"25\u00d723=575\", \"67\u00d722=1474\"],\n  [\"95\u00d716=1520\", \"59\u00d784=4956\"],\n  [\"49\u00d722=1078\", \"14\u00d789=1246\"],\n  [\"56\u00d740=2240\", \"59\u00d761=3599\"],\n  [\"41\u00d796=3936\", \"72\u00d782=5904\"],\n  [\"36\u00d717=612\", \"35\u00d778=2730\"],\n  [\"30\u00d745=1350\", \"86\u00d731=2666\"],\n  [\"95\u00d720=1900\", \"41\u00d719=779\"],\n  [\"63\u00d767=4221\", \"79\u00d732=2528\"],\n  [\"17\u00d778=1326\", \"32\u00d781=2592\"],\n  [\"44\u00d727=1188\", \"29\u00d760=1740\"],\n  [\"35\u00d756=1960\", \"20\u00d743=860\"],\n  [\"58\u00d741=2378\", \"58\u00d744=2552\"],\n  [\"61\u00d789=5429\", \"68\u00d714=952\"],\n  [\"77\u00d790=6930\", \"14\u00d736=504\"],\n  [\"17\u00d784=1428\", \"22\u00d718=396\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"25\u00d750=1250\", \"21\u00d750=1050\"),\n    @(\"94\u00d739=3666\", \"44\u00d718=792\"),\n    @(\"39\u00d753=2067\", \"64\u00d715=960\"),\n    @(\"72\u00d781=5832\", \"61\u00d794=5734\"),\n    @(\"12\u00d752=624\", \"44\u00d733=1452\"),\n    @(\"45\u00d741=1845\", \"12\u00d750=600\"),\n    @(\"64\u00d748=3072\", \"61\u00d743=2623\"),\n    @(\"22\u00d769=1518\", \"45\u00d738=1710\"),\n    @(\"76\u00d743=3268\", \"75\u00d788=6600\"),\n    @(\"25\u00d723=575\", \"67\u00d722=1474\"),\n    @(\"95\u00d716=1520\", \"59\u00d784=4956\"),\n    @(\"49\u00d722=1078\", \"14\u00d789=1246\"),\n    @(\"56\u00d740=2240\", \"59\u00d761=3599\"),\n    @(\"41\u00d796=3936\", \"72\u00d782=5904\"),\n    @(\"36\u00d717=612\", \"35\u00d778=2730\"),\n    @(\"30\u00d745=1350\", \"86\u00d731=2666\"),\n    @(\"95\u00d720=1900\", \"41\u00d719=779\"),\n    @(\"63\u00d767=4221\", \"79\u00d732=2528\"),\n    @(\"17\u00d778=1326\", \"32\u00d781=2592\"),\n    @(\"44\u00d727=1188\", \"29\u00d760=1740\"),\n    @(\"35\u00d756=1960\", \"20\u00d743=860\"),\n    @(\"58\u00d741=2378\", \"58\u00d744=2552\"),\n    @(\"61\u00d789=5429\", \"68\u00d714=952\"),\n    @(\"77\u00d790=6930\", \"14\u00d736=504\"),\n    @(\"17\u00d784=1428\", \"22\u00d718=396\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue = 1 (Wrap), wdReplaceAll = 2 (Replace)\n    $ok = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"No match found for `\"$oldText`\"\"\n    }\n}\n"}
